$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.183
$ws.Range("D5").Value = -8.148999999999999
$ws.Range("E7").Value = 13.078
$ws.Range("D9").Value = -7.746
$ws.Range("D11").Value = -8.171000000000001
$ws.Range("E11").Value = 12.949
$ws.Range("A21").Value = -20.62
$ws.Range("D21").Value = -7.936
$ws.Range("E21").Value = 13.232
$ws.Range("A23").Value = -21.368
$ws.Range("A25").Value = -21.964
